# Natmi re-run following Dr Hou's advice: the ligand/receptor-cluster
# summary table now uses 3 replicates (vs. 1) and an extra "M2" target
# cluster, growing the data block from 12 rows (2:13) to 16 rows (2:17).
# Rewrite every data row (A2:T17) with the refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Sema7a","Itgb1","ECs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,98.946724,296.840172,0.2098009692989996,0.2098009692989996,1955.417779124679,17598.76001212211,0.127811142217791,0.127811142217791),
  @("ECs","Sema7a","Itgb1","FAPs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,163.0062356666667,489.018707,0.345629090707923,0.3456290907079231,3221.382966967025,28992.44670270323,0.2105578873857319,0.2105578873857319),
  @("ECs","Sema7a","Itgb1","M2",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,65.39610666666668,196.18832,0.1386621609326595,0.1386621609326595,1292.379418863165,11631.41476976848,0.08447324733725561,0.08447324733725561),
  @("ECs","Sema7a","Itgb1","sCs",3,1,19.76232966666667,59.28698900000001,0.6092018671069146,0.6092018671069145,3,1,144.2727966666667,432.81839,0.3059077790604178,0.3059077790604179,2851.166569658635,25660.49912692771,0.186359590166136,0.186359590166136),
  @("FAPs","Sema7a","Itgb1","ECs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,98.946724,296.840172,0.2098009692989996,0.2098009692989996,968.8450936063334,8719.605842457,0.06332621057652389,0.06332621057652389),
  @("FAPs","Sema7a","Itgb1","FAPs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,163.0062356666667,489.018707,0.345629090707923,0.3456290907079231,1596.089140383139,14364.80226344825,0.1043244969395229,0.1043244969395229),
  @("FAPs","Sema7a","Itgb1","M2",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,65.39610666666668,196.18832,0.1386621609326595,0.1386621609326595,640.3314281022224,5762.98285292,0.04185371131295011,0.04185371131295011),
  @("FAPs","Sema7a","Itgb1","sCs",3,1,9.791583333333334,29.37475,0.3018394566436632,0.3018394566436631,3,1,144.2727966666667,432.81839,0.3059077790604178,0.3059077790604179,1412.659111294722,12713.9320016525,0.09233503781466629,0.09233503781466629),
  @("M2","Sema7a","Itgb1","ECs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,98.946724,296.840172,0.2098009692989996,0.2098009692989996,5.251894216471999,47.26704794824799,0.0003432773322306472,0.0003432773322306472),
  @("M2","Sema7a","Itgb1","FAPs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,163.0062356666667,489.018707,0.345629090707923,0.3456290907079231,8.652044976715333,77.86840479043799,0.0005655199430009781,0.0005655199430009781),
  @("M2","Sema7a","Itgb1","M2",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,65.39610666666668,196.18832,0.1386621609326595,0.1386621609326595,3.471094549653334,31.23985094688,0.0002268796795617426,0.0002268796795617426),
  @("M2","Sema7a","Itgb1","sCs",1,0.3333333333333333,0.05307799999999999,0.159234,0.001636204700948844,0.001636204700948844,3,1,144.2727966666667,432.81839,0.3059077790604178,0.3059077790604179,7.657711501473332,68.91940351326,0.000500527746155476,0.000500527746155476),
  @("sCs","Sema7a","Itgb1","ECs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,98.946724,296.840172,0.2098009692989996,0.2098009692989996,280.2879022579014,2522.591120321112,0.01832033917245404,0.01832033917245404),
  @("sCs","Sema7a","Itgb1","FAPs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,163.0062356666667,489.018707,0.345629090707923,0.3456290907079231,461.7502632019135,4155.752368817222,0.03018118643966735,0.03018118643966735),
  @("sCs","Sema7a","Itgb1","M2",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,65.39610666666668,196.18832,0.1386621609326595,0.1386621609326595,185.2485540949689,1667.23698685472,0.012108322602892,0.012108322602892),
  @("sCs","Sema7a","Itgb1","sCs",3,1,2.832715333333333,8.498146,0.08732247154847343,0.08732247154847342,3,1,144.2727966666667,432.81839,0.3059077790604178,0.3059077790604179,408.6837633005489,3678.15386970494,0.02671262333346003,0.02671262333346003),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $rowVals = $data[$i]
  $r = $startRow + $i
  $arr = New-Object "object[,]" 1,20
  for ($c = 0; $c -lt 20; $c++) {
    $arr[0,$c] = $rowVals[$c]
  }
  $rng = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,20))
  $rng.Value = $arr
}
